$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-27 05:03:27"
$wsZhCn.Range("H2").Value = "2016-08-27 05:03:23"
$wsZhCn.Range("K2").Value = "2016-08-27 05:03:40"
$wsDeDe.Range("K2").Value = "2016-08-27 05:03:47"
